$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: a new "Rus" column is inserted between "Date" and "Keys".
#     (Cells were retyped in place rather than via a real column insert -
#     the column-E width band stays anchored on column E.)
$ws.Range("E5").Value = "Rus"
$ws.Range("F5").Value = "Date"
$ws.Range("G5").Value = "Keys"

# --- Row 6: move the old "2024-01-15/" (date-formatted) cell from E6 to
#     F6, keep its style, then give E6 the new "get_rus_data/" content
#     with default formatting, and G6 the "KRW-USD" key.
$ws.Range("F6").Value = $ws.Range("E6").Value()
$ws.Range("E6").Copy()
$ws.Range("F6").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("E6").ClearFormats()
$ws.Range("E6").Value = "get_rus_data/"
$ws.Range("G6").Value = "KRW-USD"

# --- Row 8: URL concatenation now skips the new Rus column
$ws.Range("B8").Formula = "=CONCAT(B6,C6,F6,G6)"

# --- Row 9: the Rus endpoint URL - now a literal value instead of a formula
$ws.Range("B9").Value = "http://127.0.0.1:5000/get_rus_data/2024-01-26/Rate-USD"

# --- Row 10: WEBSERVICE + parsing/formatting pipeline for the main KRW-USD rate
$ws.Range("B10").Formula = "=WEBSERVICE(B8)"
$ws.Range("C10").Formula = "=TYPE(B10)"
$ws.Range("D10").Formula = "=SUBSTITUTE(B10, CHAR(10), """")"
$ws.Range("E10").Formula = "=TEXT(D10, ""0.00"")"
$ws.Range("F10").Formula = "=TYPE(E10)"
$ws.Range("G10").Formula = "=E10*1"

# --- Row 11: same pipeline for the Rus rate
$ws.Range("B11").Formula = "=WEBSERVICE(B9)"
$ws.Range("C11").Formula = "=TYPE(B11)"
$ws.Range("D11").Formula = "=SUBSTITUTE(B11, CHAR(10), """")"
$ws.Range("E11").Formula = "=TEXT(D11, ""0.0000"")"
$ws.Range("F11").Formula = "=TYPE(E11)"
$ws.Range("G11").Formula = "=E11*1"

# --- Row 15/16: manual example, Rus rate, quote-prefixed literal URL text
#     (written before row 13/14 - matches the order the strings were typed)
$ws.Range("B15").Value = "'http://127.0.0.1:5000/get_rus_data/2024-01-15/Rate-USD"
$ws.Range("B16").Formula = "=TEXT(VALUE(SUBSTITUTE(WEBSERVICE(B15), CHAR(10), """")), ""0.0000"")"

# --- Row 13/14: manual example, KRW-USD, quote-prefixed literal URL text
$ws.Range("B13").Value = "'http://127.0.0.1:5000/get_main_data/2024-01-15/KRW-USD"
$ws.Range("B14").Formula = "=TEXT(VALUE(SUBSTITUTE(WEBSERVICE(B13), CHAR(10), """")), ""0.00"")"

# --- Selection moved down to the next empty input row
$ws.Range("B17").Select()
